$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1:L1)
$ws.Range("I1").Value = "Invoice start date"
$ws.Range("J1").Value = "first_Price_increment_applicable_after"
$ws.Range("K1").Value = "renewal term"
$ws.Range("L1").Value = "price_increase_percentage"

# Column I holds date-like text (e.g. "2024-01-01") that must stay literal
# text rather than being auto-converted to a date serial number.
$ws.Range("I2:I56").NumberFormat = "@"

# Repeating data cycles observed in the source data
$dates = @("2024-01-01", "2024-03-15", "2024-06-01", "2024-09-01", "2025-01-01", "2025-04-01", "2025-07-01", "2025-10-01")
$incrementAfter = @(5, 6, 7)
$renewalTerm = @(12, 24)
$increasePct = @(2, 3, 4, 5)

for ($row = 2; $row -le 56; $row++) {
    $idx = $row - 2
    $ws.Cells.Item($row, 9).Value = $dates[$idx % 8]
    $ws.Cells.Item($row, 10).Value = $incrementAfter[$idx % 3]
    $ws.Cells.Item($row, 11).Value = $renewalTerm[$idx % 2]
    $ws.Cells.Item($row, 12).Value = $increasePct[$idx % 4]
}
